$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L/M inputs for rows 4-8 (koef/Unap -> Uo becomes 15 / 1.3 instead of 12 / 1.2)
$ws.Range("L4").Value = 15
$ws.Range("M4").Value = 1.3

$ws.Range("L5").Value = 15
$ws.Range("M5").Value = 1.3

$ws.Range("L6").Value = 15
$ws.Range("M6").Value = 1.3

$ws.Range("L7").Value = 15
$ws.Range("M7").Value = 1.3

$ws.Range("L8").Value = 15
$ws.Range("M8").Value = 1.3

# Add new "Rp skut" style computations for the remaining three rows (semiconductors / PV panels)
$ws.Range("T8").Formula = "=P5*Q5/(P5+Q5)"
$ws.Range("T9").Formula = "=P6*Q6/(P6+Q6)"
$ws.Range("T10").Formula = "=P7*Q7/(P7+Q7)"

# Update the selected cell shown when the sheet is reopened
$ws.Range("U17").Select()
